# Auto-generated edit script applying the cryptos.xlsx diff
# (updated price / volume figures, plus two row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.805.51"
$ws.Cells.Item(2, 5).Value = "  +0.75%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.890.05"
$ws.Cells.Item(3, 5).Value = "  +1.42%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.18%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "239.95"
$ws.Cells.Item(5, 5).Value = "  +2.37%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.0000"
$ws.Cells.Item(6, 5).Value = "  +0.16%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4773"
$ws.Cells.Item(7, 5).Value = "  +1.58%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2979"
$ws.Cells.Item(8, 5).Value = "  +8.12%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06654"
$ws.Cells.Item(9, 5).Value = "  +4.74%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "18.77"
$ws.Cells.Item(10, 5).Value = "  +7.32%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "101.54"
$ws.Cells.Item(11, 5).Value = "  +19.83%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.889.55"
$ws.Cells.Item(12, 5).Value = "  +1.54%  "

$ws.Cells.Item(13, 5).Value = "  +1.67%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.163"
$ws.Cells.Item(14, 5).Value = "  +3.91%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6648"
$ws.Cells.Item(15, 5).Value = "  +5.33%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "306.12"
$ws.Cells.Item(16, 5).Value = "  +30.05%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "30.795.17"
$ws.Cells.Item(17, 5).Value = "  +0.96%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "13.19"
$ws.Cells.Item(18, 5).Value = "  +3.98%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.000"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000007619"
$ws.Cells.Item(20, 5).Value = "  +3.65%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "2.134.80"
$ws.Cells.Item(21, 5).Value = "  +2.20%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.001"
$ws.Cells.Item(22, 5).Value = "  +0.19%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.162"
$ws.Cells.Item(23, 5).Value = "  +3.82%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.223"
$ws.Cells.Item(24, 5).Value = "  +4.71%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.366"
$ws.Cells.Item(25, 5).Value = "  +1.23%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "167.81"
$ws.Cells.Item(26, 5).Value = "  +0.43%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.73"
$ws.Cells.Item(27, 5).Value = "  +14.25%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.960"
$ws.Cells.Item(28, 5).Value = "  +3.93%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.1117"
$ws.Cells.Item(29, 5).Value = "  +8.51%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.352"
$ws.Cells.Item(30, 5).Value = "  -1.66%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.180"
$ws.Cells.Item(31, 5).Value = "  +1.87%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.007"
$ws.Cells.Item(32, 5).Value = "  +3.85%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05101"
$ws.Cells.Item(33, 5).Value = "  +3.77%  "

$ws.Cells.Item(34, 5).Value = "  +1.78%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7356"
$ws.Cells.Item(35, 5).Value = "  +3.86%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.717"
$ws.Cells.Item(36, 5).Value = "  +0.50%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01977"
$ws.Cells.Item(37, 5).Value = "  +3.10%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.709"
$ws.Cells.Item(38, 5).Value = "  +1.00%  "

$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.072"
$ws.Cells.Item(39, 5).Value = "  +5.41%  "

$ws.Cells.Item(40, 2).Value = "TrustWalletToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9070"
$ws.Cells.Item(40, 5).Value = "  +3.11%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "108.34"
$ws.Cells.Item(41, 5).Value = "  +2.40%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.9993"
$ws.Cells.Item(42, 5).Value = "  +0.05%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4214"
$ws.Cells.Item(43, 5).Value = "  +3.15%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.662"
$ws.Cells.Item(44, 5).Value = "  +2.24%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "66.97"
$ws.Cells.Item(45, 5).Value = "  +8.10%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "7.394"
$ws.Cells.Item(46, 5).Value = "  +2.36%  "

$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.107"
$ws.Cells.Item(47, 5).Value = "  +6.08%  "

$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.1234"
$ws.Cells.Item(48, 5).Value = "  -0.41%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "34.94"
$ws.Cells.Item(49, 5).Value = "  +4.04%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.05647"
$ws.Cells.Item(50, 5).Value = "  +1.63%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.403"
$ws.Cells.Item(51, 5).Value = "  +2.10%  "
